$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 932
$ws.Range("F4").Value = 47
$ws.Range("F7").Value = 1164
$ws.Range("F8").Value = 926
$ws.Range("F9").Value = 25
$ws.Range("F10").Value = 726
$ws.Range("F11").Value = 1036
$ws.Range("F12").Value = 1469
$ws.Range("F13").Value = 61
$ws.Range("F15").Value = 1625
$ws.Range("F21").Value = 1086
$ws.Range("F22").Value = 1515
$ws.Range("F24").Value = 630
$ws.Range("F25").Value = 501
$ws.Range("F26").Value = 474
$ws.Range("F30").Value = 314
$ws.Range("F33").Value = 1376
$ws.Range("F36").Value = 3978

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 179
$ws.Range("F14").Value = 4136
$ws.Range("F16").Value = 19
$ws.Range("F17").Value = 28
$ws.Range("F31").Value = 2

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 97
$ws.Range("F4").Value = 1273
$ws.Range("F5").Value = 1671
$ws.Range("F6").Value = 454
$ws.Range("F7").Value = 1013

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 97
$ws.Range("F4").Value = 1273
$ws.Range("F5").Value = 1671
$ws.Range("F6").Value = 454
$ws.Range("F7").Value = 1013
$ws.Range("F9").Value = 932
$ws.Range("F10").Value = 47
$ws.Range("F13").Value = 1164
$ws.Range("F14").Value = 926
$ws.Range("F16").Value = 25
$ws.Range("F18").Value = 726
$ws.Range("F19").Value = 179
$ws.Range("F20").Value = 179
$ws.Range("F22").Value = 1036
$ws.Range("F23").Value = 1469
$ws.Range("F24").Value = 61
$ws.Range("F26").Value = 1625
$ws.Range("F31").Value = 1086
$ws.Range("F32").Value = 1515
$ws.Range("F34").Value = 630
$ws.Range("F35").Value = 501
$ws.Range("F36").Value = 474
$ws.Range("F42").Value = 314
$ws.Range("F48").Value = 1376
$ws.Range("F50").Value = 2
$ws.Range("F51").Value = 3978
